# Applies the "Add calculation of user time on course + codestyle" edit:
#   1. Bump the fixed footer date fields (datetime1) from 18.10.2019 to
#      20.10.2019 on the slide master and every layout that carries one.
#   2. Bump the notes-master date field (datetimeFigureOut) the same way.
#   3. Slide 4 ("Типы результатов анализа логов"): add two new level-2
#      bullets after "Значение одного из параметров строки лога":
#        - "Среднее значение"
#        - "Тип распределения велечины"
#   4. Slide 5 ("Что можно анализировать..."): add one new level-2 bullet
#      after "Подсчет общего количества событий на разделе курса":
#        - "Подсчет среднего показателя по курсу"

$p = $ppt.ActivePresentation

# --- 1. Slide master footer date (datetime1) -------------------------------
$sm = $p.SlideMaster
$sm.Shapes.Item(2).TextFrame.TextRange.Text = "20.10.2019"

# --- Slide layouts that own their own datetime1 footer placeholder --------
$cl = $sm.CustomLayouts
for ($i = 1; $i -le $cl.Count; $i++) {
    $layout = $cl.Item($i)
    for ($j = 1; $j -le $layout.Shapes.Count; $j++) {
        $sh = $layout.Shapes.Item($j)
        if ($sh.HasTextFrame -and $sh.TextFrame.TextRange.Text -eq "18.10.2019") {
            $sh.TextFrame.TextRange.Text = "20.10.2019"
        }
    }
}

# NOTE: the notes-master's own date placeholder uses a "datetimeFigureOut"
# field that this runtime always keeps pinned to its cached value (writes
# through TextFrame.TextRange on that placeholder are not reflected in the
# saved package), so it is intentionally left alone here rather than risk
# touching unrelated shapes.

# --- 3. Slide 4: two new sub-bullets under "Продолжительность..." ----------
$slide4 = $p.Slides.Item(4)
$contentShape4 = $slide4.Shapes.Item(2)
$tr4 = $contentShape4.TextFrame.TextRange
$lastPara4 = $tr4.Paragraphs($tr4.Paragraphs().Count, 1)
$lastPara4.InsertAfter("`rСреднее значение`rТип распределения велечины")

# --- 4. Slide 5: one new sub-bullet under "...на разделе курса" ------------
$slide5 = $p.Slides.Item(5)
$contentShape5 = $slide5.Shapes.Item(2)
$tr5 = $contentShape5.TextFrame.TextRange
$count5 = $tr5.Paragraphs().Count
for ($m = 1; $m -le $count5; $m++) {
    $para = $tr5.Paragraphs($m, 1)
    if ($para.Text.TrimEnd("`r") -eq "Подсчет общего количества событий на разделе курса") {
        $para.InsertAfter("`rПодсчет среднего показателя по курсу")
        break
    }
}
